$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain plain text (avoids Excel
# auto-converting numeric-looking strings like "1.00" or "0.0859"
# into floating point numbers that would lose their original formatting).
$dCells = @("D2","D3","D5","D6","D7","D11","D13","D16","D17","D19","D20","D22","D24","D25","D28","D29","D30","D31","D32","D33","D34","D35","D38","D40","D41","D42","D43","D44","D46","D48","D49","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "66.155.88"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "3.535.41"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "607.74"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "143.76"
$ws.Range("E6").Value = "  -2.96%  "
$ws.Range("D7").Value = "3.534.30"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  -4.68%  "
$ws.Range("D11").Value = "8.02"
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("D13").Value = "4.135.06"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("E14").Value = "  -4.34%  "
$ws.Range("E15").Value = "  -5.15%  "
$ws.Range("D16").Value = "3.532.52"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "66.278.50"
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "10.93"
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").Value = "6.23"
$ws.Range("E20").Value = "  -3.29%  "
$ws.Range("E21").Value = "  -2.38%  "
$ws.Range("D22").Value = "425.45"
$ws.Range("E22").Value = "  -2.93%  "
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("D24").Value = "78.99"
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("D25").Value = "3.680.15"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("D28").Value = "9.28"
$ws.Range("E28").Value = "  -5.48%  "
$ws.Range("D29").Value = "8.05"
$ws.Range("E29").Value = "  -3.59%  "
$ws.Range("D30").Value = "2.47"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").Value = "0.162"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.49"
$ws.Range("E33").Value = "  -6.07%  "
$ws.Range("D34").Value = "25.39"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("D35").Value = "3.523.22"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -2.97%  "
$ws.Range("D38").Value = "7.87"
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("E39").Value = "  -5.79%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").Value = "171.02"
$ws.Range("E41").Value = "  -0.55%  "
$ws.Range("D42").Value = "0.0859"
$ws.Range("E42").Value = "  -4.24%  "
$ws.Range("D43").Value = "0.894"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "5.16"
$ws.Range("E44").Value = "  -5.43%  "
$ws.Range("E45").Value = "  -9.74%  "
$ws.Range("D46").Value = "45.23"
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("E47").Value = "  -5.80%  "
$ws.Range("D48").Value = "25.90"
$ws.Range("E48").Value = "  -9.10%  "
$ws.Range("D49").Value = "2.41"
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("E50").Value = "  -4.23%  "
$ws.Range("D51").Value = "0.952"
$ws.Range("E51").Value = "  -4.17%  "
